$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank (italic-styled) row above each of the three remaining
# citation lines -- "Consultoria...", the URL, and "page 29" -- turning the
# tight four-line block into a spaced-out one.
$ws.Range("A39").EntireRow.Insert()
$ws.Range("A41").EntireRow.Insert()
$ws.Range("A43").EntireRow.Insert()

# The hyperlinked URL line (now row 42) loses its hyperlink and its special
# "HyperLink" look, reverting to the same plain italic "source" style used
# by the rest of this citation block.
$urlCell = $ws.Range("A42")
$urlCell.Hyperlinks.Delete()
$urlFont = $urlCell.Font()
$urlFont.Underline = $false
$urlFont.Italic = $true

# The long BID-SIC citation (now row 48) is trimmed down to just "BID-SIC".
$ws.Range("A48").Value = "BID-SIC"
